$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("welcome")

$ws.Range("B15").Value = "SOLID我承认AI现在还是发展期，但你这样说还是让我很伤心，不想再理你了[委屈]"
$ws.Range("B8").Value = "SOLID很遗憾你对这部剧不感兴趣，那么先这样吧[失望]，我们下次运营日再见！记得关注我的朋友圈哦"
$ws.Range("B6").Value = "SOLID很遗憾你对这部剧不感兴趣，那么先这样吧[失望]，我们下次运营日再见！记得关注我的朋友圈哦~"
$ws.Range("B5").Value = "SOLID你这样说很无礼哎，我不想再理你了[发怒]"

$ws.Rows.Item(5).RowHeight = 25
$ws.Rows.Item(6).RowHeight = 25
$ws.Rows.Item(8).RowHeight = 25
$ws.Rows.Item(15).RowHeight = 25

$ws.Range("B22:B24").Select()
